$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new columns (F: Difference, G: Normalized Difference) ---
# Column widths for the two new columns
$ws.Columns("F").ColumnWidth = 10
$ws.Columns("G").ColumnWidth = 18

# Header labels for each of the 5 tables on the sheet
$headerRows = 2, 17, 22, 27, 32
foreach ($hr in $headerRows) {
    $ws.Range("F$hr").Value = "Difference"
    $ws.Range("G$hr").Value = "Normalized Difference"
}

# Data rows for each table: (first data row, last data row)
$dataBlocks = @(
    @(3, 14),
    @(18, 20),
    @(23, 25),
    @(28, 30),
    @(33, 35)
)

foreach ($blk in $dataBlocks) {
    $first = $blk[0]
    $last = $blk[1]
    for ($r = $first; $r -le $last; $r++) {
        $ws.Range("F$r").Formula = "=E$r-D$r"
        $ws.Range("G$r").Formula = "=F$r/E$r"
    }
}

# The second table (rows 18:20) was missing the "Numbers" (B) column values -
# fill them in (matches the 100-numbers block above it).
$ws.Range("B18").Value = 100
$ws.Range("B19").Value = 100
$ws.Range("B20").Value = 100

# --- Reposition the four existing charts to make room for the new one ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 495.633690637303
$co.Top = 34.6875590551181
$co.Width = 402.25
$co.Height = 216.375

$co = $ws.ChartObjects().Item(2)
$co.Left = 900
$co.Top = 35.4375590551181
$co.Width = 384.75
$co.Height = 216.375

$co = $ws.ChartObjects().Item(3)
$co.Left = 495.633690637303
$co.Top = 273.5625
$co.Width = 402.25
$co.Height = 216.375

$co = $ws.ChartObjects().Item(4)
$co.Left = 900
$co.Top = 276.5625
$co.Width = 384.75
$co.Height = 216.375

# --- Add the new "Normalized Expectation Error" chart ---
$newChartObj = $ws.ChartObjects().Add(495.633690637303, 539.625, 660, 387)
$chart = $newChartObj.Chart
$chart.ChartType = 51

$ser1 = $chart.SeriesCollection().NewSeries()
$ser1.Name = "=Sheet1!`$B`$3"
$ser1.XValues = "=Sheet1!`$C`$18:`$C`$20"
$ser1.Values = "=Sheet1!`$G`$3:`$G`$5"

$ser2 = $chart.SeriesCollection().NewSeries()
$ser2.Name = "=Sheet1!`$B`$6"
$ser2.XValues = "=Sheet1!`$C`$18:`$C`$20"
$ser2.Values = "=Sheet1!`$G`$6:`$G`$8"

$ser3 = $chart.SeriesCollection().NewSeries()
$ser3.Name = "=Sheet1!`$B`$9"
$ser3.XValues = "=Sheet1!`$C`$18:`$C`$20"
$ser3.Values = "=Sheet1!`$G`$9:`$G`$11"

$ser4 = $chart.SeriesCollection().NewSeries()
$ser4.Name = "=Sheet1!`$B`$12"
$ser4.XValues = "=Sheet1!`$C`$18:`$C`$20"
$ser4.Values = "=Sheet1!`$G`$12:`$G`$14"

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Normalized Expectation Error"
$chart.Legend.Position = -4160

$ws.Range("A1").Select()
